## "Error Calculations and Plots"
## Row 2 ("H 72") is removed from the missing-data table; every subsequent
## row shifts up by one. Column D/E keep the same underlying ground-truth
## values per sample, but a different subset of cells is blanked out
## (a different random "seed" of removed data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely ("H 72") - everything below shifts up by one row,
# and the sheet dimension naturally becomes A1:F62.
$ws.Rows("2:2").Delete()

# After the shift, columns A, B, C and F already line up with the target
# state. Only D and E need their "missing" markers moved to match the new
# seed's pattern - some previously-blank cells now get their value back,
# and some previously-filled cells become blank instead.

# --- Column D (col 4) fixes ---
$ws.Cells.Item(2, 4).Value = -14.8
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(20, 4).Value = -15.3
$ws.Cells.Item(21, 4).Value = ""
$ws.Cells.Item(22, 4).Value = -15.2
$ws.Cells.Item(23, 4).Value = ""
$ws.Cells.Item(24, 4).Value = -14
$ws.Cells.Item(25, 4).Value = ""
$ws.Cells.Item(52, 4).Value = -13.8
$ws.Cells.Item(53, 4).Value = ""
$ws.Cells.Item(56, 4).Value = -14.7
$ws.Cells.Item(57, 4).Value = ""
$ws.Cells.Item(58, 4).Value = -13
$ws.Cells.Item(59, 4).Value = ""

# --- Column E (col 5) fixes ---
$ws.Cells.Item(6, 5).Value = -6.4
$ws.Cells.Item(8, 5).Value = ""
$ws.Cells.Item(16, 5).Value = -5.3
$ws.Cells.Item(18, 5).Value = ""
$ws.Cells.Item(20, 5).Value = -6.9
$ws.Cells.Item(22, 5).Value = ""
$ws.Cells.Item(30, 5).Value = -7.4
$ws.Cells.Item(32, 5).Value = ""
